$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 117.057129
$ws.Range("H2").Value = 351.171387
$ws.Range("I2").Value = 0.1774070466701874
$ws.Range("J2").Value = 0.1774070466701874
$ws.Range("K2").Value = 1.0
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06674100000000001
$ws.Range("N2").Value = 0.200223
$ws.Range("O2").Value = 0.001121358778383549
$ws.Range("P2").Value = 0.001121358778383549
$ws.Range("Q2").Value = 7.812509846589
$ws.Range("R2").Value = 70.312588619301
$ws.Range("S2").Value = 0.0001989369491307146
$ws.Range("T2").Value = 0.0001989369491307146

$ws.Range("G3").Value = 117.057129
$ws.Range("H3").Value = 351.171387
$ws.Range("I3").Value = 0.1774070466701874
$ws.Range("J3").Value = 0.1774070466701874
$ws.Range("O3").Value = 0.0189158439903152
$ws.Range("P3").Value = 0.01891584399031519
$ws.Range("Q3").Value = 131.786739694414
$ws.Range("R3").Value = 1186.080657249726
$ws.Range("S3").Value = 0.003355804017595833
$ws.Range("T3").Value = 0.003355804017595832

$ws.Range("G4").Value = 117.057129
$ws.Range("H4").Value = 351.171387
$ws.Range("I4").Value = 0.1774070466701874
$ws.Range("J4").Value = 0.1774070466701874
$ws.Range("M4").Value = 31.745291
$ws.Range("N4").Value = 95.235873
$ws.Range("O4").Value = 0.5333731999099544
$ws.Range("P4").Value = 0.5333731999099544
$ws.Range("Q4").Value = 3716.012623729539
$ws.Range("R4").Value = 33444.11361356585
$ws.Range("S4").Value = 0.09462416416905249
$ws.Range("T4").Value = 0.09462416416905249

$ws.Range("G5").Value = 117.057129
$ws.Range("H5").Value = 351.171387
$ws.Range("I5").Value = 0.1774070466701874
$ws.Range("J5").Value = 0.1774070466701874
$ws.Range("M5").Value = 0.07967600000000001
$ws.Range("N5").Value = 0.239028
$ws.Range("O5").Value = 0.001338688093173426
$ws.Range("P5").Value = 0.001338688093173426
$ws.Range("Q5").Value = 9.326643810204
$ws.Range("R5").Value = 83.93979429183601
$ws.Range("S5").Value = 0.0002374927010224423
$ws.Range("T5").Value = 0.0002374927010224422

$ws.Range("G6").Value = 117.057129
$ws.Range("H6").Value = 351.171387
$ws.Range("I6").Value = 0.1774070466701874
$ws.Range("J6").Value = 0.1774070466701874
$ws.Range("M6").Value = 26.500431
$ws.Range("N6").Value = 79.501293
$ws.Range("O6").Value = 0.4452509092281735
$ws.Range("P6").Value = 0.4452509092281735
$ws.Range("Q6").Value = 3102.064370122599
$ws.Range("R6").Value = 27918.57933110339
$ws.Range("S6").Value = 0.07899064883338597
$ws.Range("T6").Value = 0.07899064883338597

$ws.Range("I7").Value = 0.2555873413068611
$ws.Range("J7").Value = 0.2555873413068611
$ws.Range("K7").Value = 1.0
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.06674100000000001
$ws.Range("N7").Value = 0.200223
$ws.Range("O7").Value = 0.001121358778383549
$ws.Range("P7").Value = 0.001121358778383549
$ws.Range("Q7").Value = 11.255351228159
$ws.Range("R7").Value = 101.298161053431
$ws.Range("S7").Value = 0.0002866051088181609
$ws.Range("T7").Value = 0.0002866051088181608

$ws.Range("I8").Value = 0.2555873413068611
$ws.Range("J8").Value = 0.2555873413068611
$ws.Range("O8").Value = 0.0189158439903152
$ws.Range("P8").Value = 0.01891584399031519
$ws.Range("S8").Value = 0.004834650274060027
$ws.Range("T8").Value = 0.004834650274060026

$ws.Range("I9").Value = 0.2555873413068611
$ws.Range("J9").Value = 0.2555873413068611
$ws.Range("M9").Value = 31.745291
$ws.Range("N9").Value = 95.235873
$ws.Range("O9").Value = 0.5333731999099544
$ws.Range("P9").Value = 0.5333731999099544
$ws.Range("Q9").Value = 5353.596740311275
$ws.Range("R9").Value = 48182.37066280148
$ws.Range("S9").Value = 0.1363234380893181
$ws.Range("T9").Value = 0.1363234380893181

$ws.Range("I10").Value = 0.2555873413068611
$ws.Range("J10").Value = 0.2555873413068611
$ws.Range("M10").Value = 0.07967600000000001
$ws.Range("N10").Value = 0.239028
$ws.Range("O10").Value = 0.001338688093173426
$ws.Range("P10").Value = 0.001338688093173426
$ws.Range("Q10").Value = 13.43673850339067
$ws.Range("R10").Value = 120.930646530516
$ws.Range("S10").Value = 0.0003421517305733475
$ws.Range("T10").Value = 0.0003421517305733474

$ws.Range("I11").Value = 0.2555873413068611
$ws.Range("J11").Value = 0.2555873413068611
$ws.Range("M11").Value = 26.500431
$ws.Range("N11").Value = 79.501293
$ws.Range("O11").Value = 0.4452509092281735
$ws.Range("P11").Value = 0.4452509092281735
$ws.Range("Q11").Value = 4469.09184163547
$ws.Range("R11").Value = 40221.82657471922
$ws.Range("S11").Value = 0.1138004961040914
$ws.Range("T11").Value = 0.1138004961040914

$ws.Range("G12").Value = 162.9464366666666
$ws.Range("H12").Value = 488.83931
$ws.Range("I12").Value = 0.2469550239393286
$ws.Range("J12").Value = 0.2469550239393286
$ws.Range("K12").Value = 1.0
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.06674100000000001
$ws.Range("N12").Value = 0.200223
$ws.Range("O12").Value = 0.001121358778383549
$ws.Range("P12").Value = 0.001121358778383549
$ws.Range("Q12").Value = 10.87520812957
$ws.Range("R12").Value = 97.87687316613
$ws.Range("S12").Value = 0.0002769251839602856
$ws.Range("T12").Value = 0.0002769251839602856

$ws.Range("G13").Value = 162.9464366666666
$ws.Range("H13").Value = 488.83931
$ws.Range("I13").Value = 0.2469550239393286
$ws.Range("J13").Value = 0.2469550239393286
$ws.Range("O13").Value = 0.0189158439903152
$ws.Range("P13").Value = 0.01891584399031519
$ws.Range("Q13").Value = 183.4504213162644
$ws.Range("R13").Value = 1651.05379184638
$ws.Range("S13").Value = 0.004671362705460894
$ws.Range("T13").Value = 0.004671362705460894

$ws.Range("G14").Value = 162.9464366666666
$ws.Range("H14").Value = 488.83931
$ws.Range("I14").Value = 0.2469550239393286
$ws.Range("J14").Value = 0.2469550239393286
$ws.Range("M14").Value = 31.745291
$ws.Range("N14").Value = 95.235873
$ws.Range("O14").Value = 0.5333731999099544
$ws.Range("P14").Value = 0.5333731999099544
$ws.Range("Q14").Value = 5172.782049396402
$ws.Range("R14").Value = 46555.03844456762
$ws.Range("S14").Value = 0.1317191913523591
$ws.Range("T14").Value = 0.1317191913523591

$ws.Range("G15").Value = 162.9464366666666
$ws.Range("H15").Value = 488.83931
$ws.Range("I15").Value = 0.2469550239393286
$ws.Range("J15").Value = 0.2469550239393286
$ws.Range("M15").Value = 0.07967600000000001
$ws.Range("N15").Value = 0.239028
$ws.Range("O15").Value = 0.001338688093173426
$ws.Range("P15").Value = 0.001338688093173426
$ws.Range("Q15").Value = 12.98292028785333
$ws.Range("R15").Value = 116.84628259068
$ws.Range("S15").Value = 0.0003305957500969376
$ws.Range("T15").Value = 0.0003305957500969376

$ws.Range("G16").Value = 162.9464366666666
$ws.Range("H16").Value = 488.83931
$ws.Range("I16").Value = 0.2469550239393286
$ws.Range("J16").Value = 0.2469550239393286
$ws.Range("M16").Value = 26.500431
$ws.Range("N16").Value = 79.501293
$ws.Range("O16").Value = 0.4452509092281735
$ws.Range("P16").Value = 0.4452509092281735
$ws.Range("Q16").Value = 4318.150801580869
$ws.Range("R16").Value = 38863.35721422783
$ws.Range("S16").Value = 0.1099569489474514
$ws.Range("T16").Value = 0.1099569489474514

$ws.Range("G17").Value = 80.88728066666667
$ws.Range("H17").Value = 242.661842
$ws.Range("I17").Value = 0.1225894885586668
$ws.Range("J17").Value = 0.1225894885586668
$ws.Range("K17").Value = 1.0
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.06674100000000001
$ws.Range("N17").Value = 0.200223
$ws.Range("O17").Value = 0.001121358778383549
$ws.Range("P17").Value = 0.001121358778383549
$ws.Range("Q17").Value = 5.398497998974001
$ws.Range("R17").Value = 48.586481990766
$ws.Range("S17").Value = 0.0001374667991328107
$ws.Range("T17").Value = 0.0001374667991328106

$ws.Range("G18").Value = 80.88728066666667
$ws.Range("H18").Value = 242.661842
$ws.Range("I18").Value = 0.1225894885586668
$ws.Range("J18").Value = 0.1225894885586668
$ws.Range("O18").Value = 0.0189158439903152
$ws.Range("P18").Value = 0.01891584399031519
$ws.Range("Q18").Value = 91.06554289236846
$ws.Range("R18").Value = 819.589886031316
$ws.Range("S18").Value = 0.002318883640428271
$ws.Range("T18").Value = 0.002318883640428271

$ws.Range("G19").Value = 80.88728066666667
$ws.Range("H19").Value = 242.661842
$ws.Range("I19").Value = 0.1225894885586668
$ws.Range("J19").Value = 0.1225894885586668
$ws.Range("M19").Value = 31.745291
$ws.Range("N19").Value = 95.235873
$ws.Range("O19").Value = 0.5333731999099544
$ws.Range("P19").Value = 0.5333731999099544
$ws.Range("Q19").Value = 2567.790262962007
$ws.Range("R19").Value = 23110.11236665807
$ws.Range("S19").Value = 0.06538594778786086
$ws.Range("T19").Value = 0.06538594778786086

$ws.Range("G20").Value = 80.88728066666667
$ws.Range("H20").Value = 242.661842
$ws.Range("I20").Value = 0.1225894885586668
$ws.Range("J20").Value = 0.1225894885586668
$ws.Range("M20").Value = 0.07967600000000001
$ws.Range("N20").Value = 0.239028
$ws.Range("O20").Value = 0.001338688093173426
$ws.Range("P20").Value = 0.001338688093173426
$ws.Range("Q20").Value = 6.444774974397334
$ws.Range("R20").Value = 58.00297476957601
$ws.Range("S20").Value = 0.0001641090886817072
$ws.Range("T20").Value = 0.0001641090886817072

$ws.Range("G21").Value = 80.88728066666667
$ws.Range("H21").Value = 242.661842
$ws.Range("I21").Value = 0.1225894885586668
$ws.Range("J21").Value = 0.1225894885586668
$ws.Range("M21").Value = 26.500431
$ws.Range("N21").Value = 79.501293
$ws.Range("O21").Value = 0.4452509092281735
$ws.Range("P21").Value = 0.4452509092281735
$ws.Range("Q21").Value = 2143.547800084634
$ws.Range("R21").Value = 19291.93020076171
$ws.Range("S21").Value = 0.05458308124256318
$ws.Range("T21").Value = 0.05458308124256317

$ws.Range("G22").Value = 130.2892406666667
$ws.Range("H22").Value = 390.867722
$ws.Range("I22").Value = 0.1974610995249561
$ws.Range("J22").Value = 0.1974610995249561
$ws.Range("K22").Value = 1.0
$ws.Range("L22").Value = 0.3333333333333333
$ws.Range("M22").Value = 0.06674100000000001
$ws.Range("N22").Value = 0.200223
$ws.Range("O22").Value = 0.001121358778383549
$ws.Range("P22").Value = 0.001121358778383549
$ws.Range("Q22").Value = 8.695634211334001
$ws.Range("R22").Value = 78.260707902006
$ws.Range("S22").Value = 0.0002214247373415771
$ws.Range("T22").Value = 0.000221424737341577

$ws.Range("G23").Value = 130.2892406666667
$ws.Range("H23").Value = 390.867722
$ws.Range("I23").Value = 0.1974610995249561
$ws.Range("J23").Value = 0.1974610995249561
$ws.Range("O23").Value = 0.0189158439903152
$ws.Range("P23").Value = 0.01891584399031519
$ws.Range("Q23").Value = 146.6838832577285
$ws.Range("R23").Value = 1320.154949319556
$ws.Range("S23").Value = 0.003735143352770171
$ws.Range("T23").Value = 0.00373514335277017

$ws.Range("G24").Value = 130.2892406666667
$ws.Range("H24").Value = 390.867722
$ws.Range("I24").Value = 0.1974610995249561
$ws.Range("J24").Value = 0.1974610995249561
$ws.Range("M24").Value = 31.745291
$ws.Range("N24").Value = 95.235873
$ws.Range("O24").Value = 0.5333731999099544
$ws.Range("P24").Value = 0.5333731999099544
$ws.Range("Q24").Value = 4136.069859132367
$ws.Range("R24").Value = 37224.6287321913
$ws.Range("S24").Value = 0.1053204585113638
$ws.Range("T24").Value = 0.1053204585113638

$ws.Range("G25").Value = 130.2892406666667
$ws.Range("H25").Value = 390.867722
$ws.Range("I25").Value = 0.1974610995249561
$ws.Range("J25").Value = 0.1974610995249561
$ws.Range("M25").Value = 0.07967600000000001
$ws.Range("N25").Value = 0.239028
$ws.Range("O25").Value = 0.001338688093173426
$ws.Range("P25").Value = 0.001338688093173426
$ws.Range("Q25").Value = 10.38092553935734
$ws.Range("R25").Value = 93.42832985421602
$ws.Range("S25").Value = 0.0002643388227989916
$ws.Range("T25").Value = 0.0002643388227989915

$ws.Range("G26").Value = 130.2892406666667
$ws.Range("H26").Value = 390.867722
$ws.Range("I26").Value = 0.1974610995249561
$ws.Range("J26").Value = 0.1974610995249561
$ws.Range("M26").Value = 26.500431
$ws.Range("N26").Value = 79.501293
$ws.Range("O26").Value = 0.4452509092281735
$ws.Range("P26").Value = 0.4452509092281735
$ws.Range("Q26").Value = 3452.721032329394
$ws.Range("R26").Value = 31074.48929096455
$ws.Range("S26").Value = 0.08791973410068155
$ws.Range("T26").Value = 0.08791973410068153
